# Update BOC USD rates (auto)
# Adds the 2026-01-03 05:30:00 publish row to "All Published Values" and
# bumps the 2026-01-03 publish count on "Daily Summary".

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("All Published Values")
$ws2 = $wb.Worksheets.Item("Daily Summary")

# Helper: write a value as literal text, even when it looks like a
# number/date, without leaving a stray NumberFormat behind on the cell.
function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# --- New row 21 on "All Published Values" -------------------------------
Set-TextValue $ws1.Range("A21") "2026-01-03"
$ws1.Range("B21").Value = "2026-01-03 05:30:00"
Set-TextValue $ws1.Range("C21") "697.85"
Set-TextValue $ws1.Range("D21") "697.85"
Set-TextValue $ws1.Range("E21") "700.79"
Set-TextValue $ws1.Range("F21") "700.79"
Set-TextValue $ws1.Range("G21") "702.88"
$ws1.Range("H21").Value = "2026/01/03 05:30:00"
$ws1.Range("I21").Value = "2026-01-02 21:37:51"
$ws1.Range("J21").Value = "https://www.bankofchina.com/sourcedb/whpj/enindex_1619.html"

# Re-apply the AutoFilter so it (and the hidden _FilterDatabase name) spans
# the newly added row.
$ws1.AutoFilterMode = $false
$ws1.Range("A1:J21").AutoFilter()

foreach ($n in $wb.Names) {
    if ($n.Name -eq "All Published Values!_FilterDatabase") {
        $n.RefersTo = "='All Published Values'!`$A`$1:`$J`$21"
    }
}

# --- "Daily Summary": publishes count for 2026-01-03 goes from 1 to 2 ---
$ws2.Range("B5").Value = 2

Write-Output "done"
